# Auto-generated edit script: applies numeric updates to the Atomos_Profits data
# (split across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 29188.693
$ws.Range("J17").Value = 29188.693
$ws.Range("L17").Value = 87566.079
$ws.Range("N17").Value = -87902.079

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 668
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 697.1429000000001
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 697.1429000000001
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1349.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 18636.557
$ws.Range("I58").Value = 372.85715
$ws.Range("J58").Value = 21004.074
$ws.Range("K58").Value = 1118.57145
$ws.Range("L58").Value = 63012.222
$ws.Range("M58").Value = -968.5714499999999
$ws.Range("N58").Value = -63312.222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4074.875
$ws.Range("I74").Value = 3749
$ws.Range("J74").Value = 4121.4287
$ws.Range("K74").Value = 3749
$ws.Range("L74").Value = 4121.4287
$ws.Range("M74").Value = -2813
$ws.Range("N74").Value = -5993.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4074.875
$ws.Range("I77").Value = 3749
$ws.Range("J77").Value = 4121.4287
$ws.Range("K77").Value = 18745
$ws.Range("L77").Value = 20607.1435
$ws.Range("M77").Value = -14065
$ws.Range("N77").Value = -29967.1435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 452.42105
$ws.Range("I103").Value = 775
$ws.Range("J103").Value = 391.9375
$ws.Range("K103").Value = 2325
$ws.Range("L103").Value = 1175.8125
$ws.Range("M103").Value = -1739
$ws.Range("N103").Value = -2347.8125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1641519.9
$ws.Range("I137").Value = 2274937.8
$ws.Range("J137").Value = 2085.2942
$ws.Range("K137").Value = 6824813.399999999
$ws.Range("L137").Value = 6255.882599999999
$ws.Range("M137").Value = -6822263.399999999
$ws.Range("N137").Value = -11355.8826

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 27477.857
$ws.Range("I139").Value = 12345
$ws.Range("K139").Value = 12345
$ws.Range("M139").Value = -7205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4138
$ws.Range("I102").Value = 3161
$ws.Range("K102").Value = 3161
$ws.Range("M102").Value = -1539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1645.3334
$ws.Range("I105").Value = 1490.909
$ws.Range("J105").Value = 1815.2
$ws.Range("K105").Value = 1490.909
$ws.Range("L105").Value = 1815.2
$ws.Range("M105").Value = 256.0909999999999
$ws.Range("N105").Value = -5309.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 19999
$ws.Range("J132").Value = 19999
$ws.Range("L132").Value = 19999
$ws.Range("N132").Value = -30119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5075.1113
$ws.Range("I99").Value = 3497
$ws.Range("J99").Value = 6337.6
$ws.Range("K99").Value = 3497
$ws.Range("L99").Value = 6337.6
$ws.Range("M99").Value = -1999
$ws.Range("N99").Value = -9333.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5075.1113
$ws.Range("I126").Value = 3497
$ws.Range("J126").Value = 6337.6
$ws.Range("K126").Value = 10491
$ws.Range("L126").Value = 19012.8
$ws.Range("M126").Value = -8021
$ws.Range("N126").Value = -23952.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2248.6453
$ws.Range("I132").Value = 1618.8148
$ws.Range("K132").Value = 4856.4444
$ws.Range("M132").Value = -2326.4444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1889.4193
$ws.Range("I134").Value = 950.3200000000001
$ws.Range("J134").Value = 5802.3335
$ws.Range("K134").Value = 2850.96
$ws.Range("L134").Value = 17407.0005
$ws.Range("M134").Value = -315.96
$ws.Range("N134").Value = -22477.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1810.1818
$ws.Range("J5").Value = 2333.3333
$ws.Range("L5").Value = 6999.999899999999
$ws.Range("N5").Value = -7223.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1792.2667
$ws.Range("I131").Value = 2423.6365
$ws.Range("J131").Value = 1426.7368
$ws.Range("K131").Value = 7270.9095
$ws.Range("L131").Value = 4280.2104
$ws.Range("M131").Value = -2230.9095
$ws.Range("N131").Value = -14360.2104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1810.1818
$ws.Range("J135").Value = 2333.3333
$ws.Range("L135").Value = 20999.9997
$ws.Range("N135").Value = -26069.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3421.9524
$ws.Range("I122").Value = 2643.25
$ws.Range("J122").Value = 4129.864
$ws.Range("K122").Value = 7929.75
$ws.Range("L122").Value = 12389.592
$ws.Range("M122").Value = -5479.75
$ws.Range("N122").Value = -17289.592

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2510.6177
$ws.Range("I132").Value = 2070.44
$ws.Range("J132").Value = 3733.3333
$ws.Range("K132").Value = 6211.32
$ws.Range("L132").Value = 11199.9999
$ws.Range("M132").Value = -3681.32
$ws.Range("N132").Value = -16259.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1291.6666
$ws.Range("I7").Value = 1187.5
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1187.5
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1075.5
$ws.Range("N7").Value = -1724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2356.75
$ws.Range("I46").Value = 533.6667
$ws.Range("J46").Value = 2964.4443
$ws.Range("K46").Value = 533.6667
$ws.Range("L46").Value = 2964.4443
$ws.Range("M46").Value = -345.6667
$ws.Range("N46").Value = -3340.4443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2224.2415
$ws.Range("I68").Value = 1068.1818
$ws.Range("J68").Value = 5857.5713
$ws.Range("K68").Value = 1068.1818
$ws.Range("L68").Value = 5857.5713
$ws.Range("M68").Value = -319.1818000000001
$ws.Range("N68").Value = -7355.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2224.2415
$ws.Range("I71").Value = 1068.1818
$ws.Range("J71").Value = 5857.5713
$ws.Range("K71").Value = 5340.909000000001
$ws.Range("L71").Value = 29287.8565
$ws.Range("M71").Value = -1596.909000000001
$ws.Range("N71").Value = -36775.85649999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2463.5293
$ws.Range("I82").Value = 1766.6666
$ws.Range("J82").Value = 4136
$ws.Range("K82").Value = 1766.6666
$ws.Range("L82").Value = 4136
$ws.Range("M82").Value = -1405.6666
$ws.Range("N82").Value = -4858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2463.5293
$ws.Range("I85").Value = 1766.6666
$ws.Range("J85").Value = 4136
$ws.Range("K85").Value = 1766.6666
$ws.Range("L85").Value = 4136
$ws.Range("M85").Value = -518.6666
$ws.Range("N85").Value = -6632

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2372.7273
$ws.Range("J100").Value = 2775
$ws.Range("L100").Value = 2775
$ws.Range("N100").Value = -3857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1291.6666
$ws.Range("I126").Value = 1187.5
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 3562.5
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -1092.5
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1937.1177
$ws.Range("I132").Value = 1180.0741
$ws.Range("J132").Value = 4857.143
$ws.Range("K132").Value = 3540.2223
$ws.Range("L132").Value = 14571.429
$ws.Range("M132").Value = -1010.2223
$ws.Range("N132").Value = -19631.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 558257.2
$ws.Range("I122").Value = 1113136.6
$ws.Range("K122").Value = 3339409.8
$ws.Range("M122").Value = -3336959.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4168399.2
$ws.Range("I126").Value = 1040.7894
$ws.Range("J126").Value = 20004362
$ws.Range("K126").Value = 3122.3682
$ws.Range("L126").Value = 60013086
$ws.Range("M126").Value = -652.3681999999999
$ws.Range("N126").Value = -60018026

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 258928.48
$ws.Range("I132").Value = 387888.84
$ws.Range("J132").Value = 35397.2
$ws.Range("K132").Value = 1163666.52
$ws.Range("L132").Value = 106191.6
$ws.Range("M132").Value = -1161136.52
$ws.Range("N132").Value = -111251.6

